$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 98: Friss mi Stoub / noREST scores were recomputed ---
$ws.Range("C98").Value = 8.477151835919745
$ws.Range("F98").Value = 8.477151835919745

# --- Append four new rows for 2025-02-26 (abs_activity, rel_activity, abs_sleep, rel_sleep) ---
# Note: Range.Value performs Excel's usual "smart" text->date coercion, so typing
# "2025-02-26" directly into column A would silently turn it into a date serial.
# Building the literal text in a scratch cell and Copy/PasteSpecial-ing it keeps it
# as plain text, matching how the other date cells in this column are stored.
$ws.Range("H1").Formula = "=""2025-02-26"""
$ws.Range("H1").Copy()
$ws.Range("A102").PasteSpecial()
$ws.Range("A103").PasteSpecial()
$ws.Range("A104").PasteSpecial()
$ws.Range("A105").PasteSpecial()
$ws.Range("H1").ClearContents()

$ws.Range("B102").Value = "abs_activity"
$ws.Range("C102").Value = 7.776348091509672
$ws.Range("D102").Value = 0
$ws.Range("E102").Value = 0
$ws.Range("F102").Value = 7.776348091509672

$ws.Range("B103").Value = "rel_activity"
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 0
$ws.Range("E103").Value = 0
$ws.Range("F103").Value = 0

$ws.Range("B104").Value = "abs_sleep"
$ws.Range("C104").Value = 10
$ws.Range("D104").Value = 0
$ws.Range("E104").Value = 0
$ws.Range("F104").Value = 10

$ws.Range("B105").Value = "rel_sleep"
$ws.Range("C105").Value = 10
$ws.Range("D105").Value = 0
$ws.Range("E105").Value = 0
$ws.Range("F105").Value = 10
